$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
Write-Output $d.Name
$d.Name = "TestDesignName"
Write-Output "set done"
